# Updates cryptos list prices/volumes (GitHub Actions scheduled refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a string into a cell while forcing it to remain a Text
# cell (the source sheet stores Price/Volume as inline strings, even when
# they look numeric, e.g. "338.03" or "0.9988"). Flip to text format only
# for the write, then restore General/Normal so the cell's style index is
# unchanged (matches the original, unstyled cells).
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
    $range.Style = "Normal"
}

# Row 40/41: Algorand and Hedera swapped rank positions, with refreshed data.
Set-TextValue $ws.Range("B40") "Hedera"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D40") "0.06394"
Set-TextValue $ws.Range("E40") "  +2.55%  "

Set-TextValue $ws.Range("B41") "Algorand"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D41") "0.2217"
Set-TextValue $ws.Range("E41") "  +1.86%  "

# Price (D) / Volume(1h) (E) refresh for every other row.
$updates = @{
    2  = @{ D = "27.432.86";    E = "  +2.51%  " }
    3  = @{ D = "1.800.60";     E = "  +3.60%  " }
    4  = @{ E = "  +0.44%  " }
    5  = @{ D = "338.03";       E = "  +1.34%  " }
    6  = @{ D = "0.9988";       E = "  +0.23%  " }
    7  = @{ D = "0.3809";       E = "  +1.90%  " }
    8  = @{ D = "0.3471";       E = "  +2.23%  " }
    9  = @{ D = "48.59";        E = "  +0.80%  " }
    10 = @{ D = "1.208";        E = "  +1.75%  " }
    11 = @{ D = "0.07532";      E = "  +1.02%  " }
    12 = @{ D = "0.9987";       E = "  +0.17%  " }
    13 = @{ D = "22.09";        E = "  +8.23%  " }
    14 = @{ D = "6.506";        E = "  +1.37%  " }
    15 = @{ D = "1.796.58";     E = "  +3.57%  " }
    16 = @{ D = "7.077";        E = "  -0.05%  " }
    17 = @{ D = "0.00001103";   E = "  +2.44%  " }
    18 = @{ D = "0.06657";      E = "  -1.08%  " }
    19 = @{ D = "85.18";        E = "  +3.05%  " }
    20 = @{ D = "0.9992";       E = "  +0.28%  " }
    21 = @{ D = "6.528";        E = "  +4.91%  " }
    22 = @{ D = "17.40";        E = "  +4.32%  " }
    23 = @{ D = "27.411.96";    E = "  +2.53%  " }
    24 = @{ D = "12.56";        E = "  -1.58%  " }
    25 = @{ D = "2.440";        E = "  -0.18%  " }
    26 = @{ D = "2.587";        E = "  +6.67%  " }
    27 = @{ D = "1.508";        E = "  +1.43%  " }
    28 = @{ D = "21.48";        E = "  +9.87%  " }
    29 = @{ D = "153.00";       E = "  +1.05%  " }
    30 = @{ D = "1.999.31";     E = "  +3.54%  " }
    31 = @{ D = "134.51";       E = "  +1.66%  " }
    32 = @{ D = "4.060";        E = "  -0.87%  " }
    33 = @{ D = "6.170";        E = "  +1.99%  " }
    34 = @{ D = "0.08718";      E = "  +0.90%  " }
    35 = @{ D = "13.33";        E = "  +3.72%  " }
    36 = @{ D = "1.694";        E = "  -0.30%  " }
    37 = @{ D = "5.484";        E = "  +1.58%  " }
    38 = @{ D = "0.6918";       E = "  +10.47%  " }
    39 = @{ D = "8.966";        E = "  +5.97%  " }
    42 = @{ D = "0.02342";      E = "  -0.33%  " }
    43 = @{ E = "  +4.12%  " }
    44 = @{ D = "14.41";        E = "  +1.39%  " }
    45 = @{ D = "0.6479";       E = "  +6.24%  " }
    46 = @{ D = "0.9988";       E = "  +0.35%  " }
    47 = @{ D = "3.874";        E = "  -1.10%  " }
    48 = @{ D = "2.137";        E = "  +3.67%  " }
    49 = @{ D = "130.28";       E = "  +0.82%  " }
    50 = @{ D = "0.07200";      E = "  -0.18%  " }
    51 = @{ D = "79.64";        E = "  +2.56%  " }
}

foreach ($row in $updates.Keys) {
    $cells = $updates[$row]
    foreach ($col in $cells.Keys) {
        $ref = "$col$row"
        Set-TextValue $ws.Range($ref) $cells[$col]
    }
}
